# Apply updated ticket/stock counts (and a couple of status-text tweaks)
# to the 广州-漫展信息 workbook, as published for commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value = 61
$wsExpo.Range("F13").Value = 99
$wsExpo.Range("F14").Value = 952
$wsExpo.Range("F16").Value = 2068
$wsExpo.Range("F18").Value = 9091
$wsExpo.Range("F20").Value = 532

# --- Sheet "演出" (Performances) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 528
$wsShow.Range("F13").Value = 23

# --- Sheet "本地生活" (Local Life) ---
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F2").Value = 5622
$wsLocal.Range("G2").Value = "暂时售罄"
$wsLocal.Range("F3").Value = 435

# --- Sheet "全部类型" (All Types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 5622
$wsAll.Range("G3").Value = "暂时售罄"
$wsAll.Range("F4").Value = 435
$wsAll.Range("F6").Value = 528
$wsAll.Range("F10").Value = 61
$wsAll.Range("F20").Value = 99
$wsAll.Range("F22").Value = 952
$wsAll.Range("F27").Value = 2068
$wsAll.Range("F29").Value = 9091
$wsAll.Range("F31").Value = 23
$wsAll.Range("F33").Value = 532
